$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","E8","D9","E9","D10","E10","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","E17","D18","E18","D19","E19","E20","D21","E21","D22","E22","D23","E23","E24","D25","E25","D26","E26","E27","E39","D40","E40","D41","E41","D42","E42","E43","D44","E44","D45","E45","D46","E46","E47","D49","E49","E50","E51")
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "303.02"
$ws.Range("E2").Value = "-0.34%"
$ws.Range("D3").Value = "37.04"
$ws.Range("E3").Value = "5.78%"
$ws.Range("D4").Value = "5.014"
$ws.Range("E4").Value = "-3.07%"
$ws.Range("D5").Value = "0.07865"
$ws.Range("E5").Value = "0.69%"
$ws.Range("D6").Value = "2.210"
$ws.Range("E6").Value = "-3.68%"
$ws.Range("D7").Value = "7.996"
$ws.Range("E7").Value = "-0.69%"
$ws.Range("E8").Value = "0.61%"
$ws.Range("D9").Value = "0.9207"
$ws.Range("E9").Value = "-0.34%"
$ws.Range("D10").Value = "0.09561"
$ws.Range("E10").Value = "-4.95%"
$ws.Range("E11").Value = "2.39%"
$ws.Range("D12").Value = "0.08573"
$ws.Range("E12").Value = "0.73%"
$ws.Range("D13").Value = "0.03592"
$ws.Range("E13").Value = "6.31%"
$ws.Range("D14").Value = "0.09974"
$ws.Range("E14").Value = "0.59%"
$ws.Range("D15").Value = "0.001497"
$ws.Range("E15").Value = "0.63%"
$ws.Range("D16").Value = "0.005696"
$ws.Range("E16").Value = "-1.77%"
$ws.Range("E17").Value = "-0.64%"
$ws.Range("D18").Value = "2.249"
$ws.Range("E18").Value = "5.75%"
$ws.Range("D19").Value = "0.3420"
$ws.Range("E19").Value = "0.09%"
$ws.Range("E20").Value = "-0.79%"
$ws.Range("D21").Value = "4.768"
$ws.Range("E21").Value = "4.55%"
$ws.Range("D22").Value = "0.2199"
$ws.Range("E22").Value = "-8.23%"
$ws.Range("D23").Value = "0.04595"
$ws.Range("E23").Value = "-1.18%"
$ws.Range("E24").Value = "0.74%"
$ws.Range("D25").Value = "0.004464"
$ws.Range("E25").Value = "3.15%"
$ws.Range("D26").Value = "0.0001398"
$ws.Range("E26").Value = "7.54%"
$ws.Range("E27").Value = "39.64%"
$ws.Range("E39").Value = "4.14%"
$ws.Range("D40").Value = "0.04726"
$ws.Range("E40").Value = "-0.31%"
$ws.Range("D41").Value = "0.008108"
$ws.Range("E41").Value = "5.52%"
$ws.Range("D42").Value = "0.1395"
$ws.Range("E42").Value = "-1.33%"
$ws.Range("E43").Value = "3.87%"
$ws.Range("D44").Value = "0.002227"
$ws.Range("E44").Value = "0.81%"
$ws.Range("D45").Value = "0.01041"
$ws.Range("E45").Value = "3.46%"
$ws.Range("D46").Value = "0.00006144"
$ws.Range("E46").Value = "1.83%"
$ws.Range("E47").Value = "-0.07%"
$ws.Range("D49").Value = "10.71"
$ws.Range("E49").Value = "176.05%"
$ws.Range("E50").Value = "-0.05%"
$ws.Range("E51").Value = "-0.07%"
